$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.183.63'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '3.556.50'
$ws.Range("E3").Value = '  +1.23%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.12'
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.53'
$ws.Range("E6").Value = '  -1.00%  '

$ws.Range("D7").Value = '3.555.90'
$ws.Range("E7").Value = '  +1.22%  '

$ws.Range("E9").Value = '  +2.47%  '

$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.83'
$ws.Range("E11").Value = '  -2.05%  '

$ws.Range("E12").Value = '  -0.41%  '

$ws.Range("D13").Value = '4.159.11'
$ws.Range("E13").Value = '  +1.36%  '

$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.01'
$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("D16").Value = '3.561.65'
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("D17").Value = '66.256.76'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.35'
$ws.Range("E19").Value = '  +7.84%  '

$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.66'
$ws.Range("E21").Value = '  -1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '428.90'
$ws.Range("E22").Value = '  +0.46%  '

$ws.Range("E23").Value = '  +1.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.72'
$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("D25").Value = '3.698.57'
$ws.Range("E25").Value = '  +1.64%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  -3.12%  '

$ws.Range("E28").Value = '  +0.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.07'
$ws.Range("E29").Value = '  -2.46%  '

$ws.Range("E30").Value = '  -1.85%  '

$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("D32").Value = '3.552.24'

$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("E34").Value = '  -1.89%  '

$ws.Range("E35").Value = '  -9.35%  '

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.80'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("E38").Value = '  -1.68%  '

$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '173.91'
$ws.Range("E40").Value = '  +2.03%  '

$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.19'
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.886'
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("E44").Value = '  +0.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '45.89'
$ws.Range("E45").Value = '  +1.19%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.81'
$ws.Range("E48").Value = '  -4.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.39'
$ws.Range("E49").Value = '  -2.60%  '

$ws.Range("E50").Value = '  -1.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.88'
$ws.Range("E51").Value = '  +1.36%  '
